# Applies the "chore: update Sheets via scheduled runner" numeric updates
# to the Kraken_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7742.7144
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 8833.166999999999
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 8833.166999999999
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -9183.166999999999

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H70").Value = 37999.332
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 37999.332
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 113997.996
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -114537.996

$ws.Range("H73").Value = 37999.332
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 37999.332
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 113997.996
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -115869.996

$ws.Range("H98").Value = 1742
$ws.Range("I98").Value = 1472.75
$ws.Range("K98").Value = 1472.75
$ws.Range("M98").Value = 25.25

$ws.Range("H110").Value = 89997
$ws.Range("J110").Value = 89997
$ws.Range("L110").Value = 89997
$ws.Range("N110").Value = -98177

$ws.Range("H122").Value = 1742
$ws.Range("I122").Value = 1472.75
$ws.Range("K122").Value = 4418.25
$ws.Range("M122").Value = -1968.25

$ws.Range("H138").Value = 3812.2354
$ws.Range("J138").Value = 3993.8667
$ws.Range("L138").Value = 11981.6001
$ws.Range("N138").Value = -22261.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 861.58826
$ws.Range("I63").Value = 926.4666999999999
$ws.Range("J63").Value = 375
$ws.Range("K63").Value = 926.4666999999999
$ws.Range("L63").Value = 375
$ws.Range("M63").Value = -240.4666999999999
$ws.Range("N63").Value = -1747

$ws.Range("H66").Value = 861.58826
$ws.Range("I66").Value = 926.4666999999999
$ws.Range("J66").Value = 375
$ws.Range("K66").Value = 4632.3335
$ws.Range("L66").Value = 1875
$ws.Range("M66").Value = -1200.3335
$ws.Range("N66").Value = -8739

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 483.33334
$ws.Range("I107").Value = 483.33334
$ws.Range("K107").Value = 483.33334
$ws.Range("M107").Value = 1436.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5105.2856
$ws.Range("J31").Value = 6347.5
$ws.Range("L31").Value = 6347.5
$ws.Range("N31").Value = -6937.5

$ws.Range("H34").Value = 5105.2856
$ws.Range("J34").Value = 6347.5
$ws.Range("L34").Value = 6347.5
$ws.Range("N34").Value = -6751.5

$ws.Range("H62").Value = 125001000
$ws.Range("I62").Value = 125001000
$ws.Range("K62").Value = 125001000
$ws.Range("M62").Value = -125000376

$ws.Range("H65").Value = 125001000
$ws.Range("I65").Value = 125001000
$ws.Range("K65").Value = 625005000
$ws.Range("M65").Value = -625001880

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2977.4443
$ws.Range("J134").Value = 3249.5
$ws.Range("L134").Value = 9748.5
$ws.Range("N134").Value = -14818.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 850
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3182

$ws.Range("H139").Value = 1589.8
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H80").Value = 69668.664
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 69668.664
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 69668.664
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -71664.664

$ws.Range("H83").Value = 69668.664
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 69668.664
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 348343.32
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -358327.32

$ws.Range("H122").Value = 1566.25
$ws.Range("I122").Value = 1378.5
$ws.Range("J122").Value = 1754
$ws.Range("K122").Value = 4135.5
$ws.Range("L122").Value = 5262
$ws.Range("M122").Value = -1685.5
$ws.Range("N122").Value = -10162

$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

$ws.Range("H140").Value = 150000
$ws.Range("J140").Value = 150000
$ws.Range("L140").Value = 150000
$ws.Range("N140").Value = -160360

$ws.Range("H141").Value = 72500
$ws.Range("J141").Value = 72500
$ws.Range("L141").Value = 72500
$ws.Range("N141").Value = -82860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3094.4595
$ws.Range("I46").Value = 2999.5
$ws.Range("J46").Value = 3099.8857
$ws.Range("K46").Value = 2999.5
$ws.Range("L46").Value = 3099.8857
$ws.Range("M46").Value = -2811.5
$ws.Range("N46").Value = -3475.8857

$ws.Range("H68").Value = 2714.2856
$ws.Range("I68").Value = 2714.2856
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2714.2856
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1965.2856
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 2714.2856
$ws.Range("I71").Value = 2714.2856
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13571.428
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9827.428
$ws.Range("N71").ClearContents()

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1788.9
$ws.Range("I122").Value = 1988.125
$ws.Range("K122").Value = 5964.375
$ws.Range("M122").Value = -3514.375

$ws.Range("H132").Value = 2666.6667
$ws.Range("I132").Value = 2666.6667
$ws.Range("K132").Value = 8000.000100000001
$ws.Range("M132").Value = -5470.000100000001
